# Apply cell updates from the crypto price refresh diff.
# Every written cell is forced to Text format ("@") before assignment so that
# numeric-looking strings (e.g. "69.00", "1.672") are preserved verbatim as
# text rather than being auto-coerced into numbers by Excel. The style is then
# reset back to "Normal" so no stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "27.933.84"
Set-TextValue "E2" "  +7.00%  "
# Row 3
Set-TextValue "D3" "1.745.60"
# Row 4
Set-TextValue "E4" "  -0.10%  "
# Row 5
Set-TextValue "D5" "228.19"
Set-TextValue "E5" "  +4.38%  "
# Row 6
Set-TextValue "D6" "0.5472"
Set-TextValue "E6" "  +4.16%  "
# Row 7
Set-TextValue "E7" "  -0.16%  "
# Row 8
Set-TextValue "E8" "  +4.21%  "
# Row 9
Set-TextValue "D9" "0.06757"
Set-TextValue "E9" "  +6.23%  "
# Row 10
Set-TextValue "D10" "21.87"
Set-TextValue "E10" "  +6.21%  "
# Row 11
Set-TextValue "D11" "0.07787"
Set-TextValue "E11" "  +1.10%  "
# Row 12
Set-TextValue "D12" "4.713"
Set-TextValue "E12" "  +2.44%  "
# Row 13
Set-TextValue "D13" "1.744.55"
Set-TextValue "E13" "  +6.42%  "
# Row 14
Set-TextValue "D14" "1.984.65"
Set-TextValue "E14" "  +5.53%  "
# Row 15
Set-TextValue "D15" "0.5996"
Set-TextValue "E15" "  +6.91%  "
# Row 16
Set-TextValue "D16" "0.0₅8422"
Set-TextValue "E16" "  +2.44%  "
# Row 17
Set-TextValue "D17" "69.00"
Set-TextValue "E17" "  +5.51%  "
# Row 18
Set-TextValue "D18" "27.941.89"
Set-TextValue "E18" "  +7.01%  "
# Row 19
Set-TextValue "D19" "225.44"
Set-TextValue "E19" "  +17.93%  "
# Row 20
Set-TextValue "D20" "4.859"
Set-TextValue "E20" "  +3.34%  "
# Row 21
Set-TextValue "E21" "  -0.17%  "
# Row 22
Set-TextValue "D22" "10.95"
Set-TextValue "E22" "  +5.36%  "
# Row 23
Set-TextValue "D23" "6.258"
Set-TextValue "E23" "  +4.54%  "
# Row 24
Set-TextValue "D24" "1.004"
Set-TextValue "E24" "  -0.10%  "
# Row 25
Set-TextValue "D25" "146.05"
Set-TextValue "E25" "  -0.01%  "
# Row 26
Set-TextValue "D26" "0.1255"
Set-TextValue "E26" "  +4.44%  "
# Row 27
Set-TextValue "B27" "Toncoin"
Set-TextValue "C27" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D27" "1.672"
Set-TextValue "E27" "  +11.77%  "
# Row 28
Set-TextValue "B28" "Cosmos"
Set-TextValue "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D28" "7.484"
Set-TextValue "E28" "  +3.08%  "
# Row 29
Set-TextValue "B29" "EthereumClassic"
Set-TextValue "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D29" "17.24"
Set-TextValue "E29" "  +8.13%  "
# Row 30
Set-TextValue "D30" "0.05691"
Set-TextValue "E30" "  +0.82%  "
# Row 32
Set-TextValue "E32" "  +5.80%  "
# Row 33
Set-TextValue "D33" "3.532"
Set-TextValue "E33" "  +4.26%  "
# Row 34
Set-TextValue "D34" "1.692"
Set-TextValue "E34" "  +6.87%  "
# Row 35
Set-TextValue "D35" "0.9784"
Set-TextValue "E35" "  +3.24%  "
# Row 36
Set-TextValue "D36" "2.858"
Set-TextValue "E36" "  +2.23%  "
# Row 37
Set-TextValue "D37" "2.447"
Set-TextValue "E37" "  +1.65%  "
# Row 38
Set-TextValue "D38" "0.5996"
Set-TextValue "E38" "  +3.60%  "
# Row 39
Set-TextValue "E39" "  +4.81%  "
# Row 40
Set-TextValue "D40" "5.943"
Set-TextValue "E40" "  -0.53%  "
# Row 41
Set-TextValue "D41" "0.8522"
Set-TextValue "E41" "  +1.42%  "
# Row 42
Set-TextValue "D42" "1.047.35"
Set-TextValue "E42" "  +2.57%  "
# Row 43
Set-TextValue "E43" "  -0.07%  "
# Row 44
Set-TextValue "D44" "102.27"
Set-TextValue "E44" "  +0.59%  "
# Row 45
Set-TextValue "D45" "1.890.34"
# Row 46
Set-TextValue "E46" "  +12.69%  "
# Row 47
Set-TextValue "D47" "59.85"
Set-TextValue "E47" "  +2.41%  "
# Row 48
Set-TextValue "D48" "8.293"
Set-TextValue "E48" "  +3.13%  "
# Row 49
Set-TextValue "D49" "0.4437"
Set-TextValue "E49" "  +2.17%  "
# Row 50
Set-TextValue "B50" "Frax"
Set-TextValue "C50" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D50" "1.003"
Set-TextValue "E50" "  -0.08%  "
# Row 51
Set-TextValue "B51" "Cronos"
Set-TextValue "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.05336"
Set-TextValue "E51" "  -0.04%  "
